$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value is a plain number; force Text format so
# Excel keeps storing them as text (matching column D's existing text cells)
# instead of silently converting to a numeric value. The style is reset back
# to Normal right after the write so the cell's formatting/style index is
# left exactly as it was before the edit.
$textCells = @("D5", "D6", "D7", "D9", "D11", "D13", "D14", "D15", "D17", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D33", "D34", "D36", "D40", "D42", "D44", "D46", "D47", "D49", "D51")

$updates = @{
    "D2" = "36.545.45"
    "E2" = "  -2.62%  "
    "D3" = "1.985.25"
    "E3" = "  -3.43%  "
    "E4" = "  +0.02%  "
    "D5" = "244.67"
    "E5" = "  -3.15%  "
    "D6" = "0.629"
    "E6" = "  -4.09%  "
    "D7" = "59.34"
    "E7" = "  -9.49%  "
    "E8" = "  +0.05%  "
    "D9" = "0.376"
    "E9" = "  -2.08%  "
    "D11" = "0.0826"
    "E11" = "  +7.31%  "
    "E12" = "  -0.94%  "
    "D13" = "23.69"
    "E13" = "  +5.07%  "
    "D14" = "0.866"
    "E14" = "  -6.03%  "
    "D15" = "14.04"
    "E15" = "  -5.53%  "
    "D16" = "2.275.99"
    "E16" = "  -3.42%  "
    "D17" = "5.47"
    "E17" = "  -2.15%  "
    "D18" = "1.988.02"
    "E18" = "  -3.09%  "
    "D19" = "36.386.06"
    "E19" = "  -2.68%  "
    "D20" = "70.44"
    "E20" = "  -4.47%  "
    "D21" = "0.0₃0865"
    "E21" = "  -1.45%  "
    "D22" = "5.34"
    "E22" = "  -2.77%  "
    "D23" = "234.14"
    "E23" = "  -2.64%  "
    "D24" = "1.00"
    "E24" = "  +0.05%  "
    "D25" = "2.60"
    "E25" = "  -1.25%  "
    "D26" = "2.30"
    "E26" = "  -4.27%  "
    "D27" = "10.09"
    "E27" = "  -0.09%  "
    "D28" = "162.19"
    "E28" = "  +0.04%  "
    "D29" = "19.87"
    "E29" = "  -0.78%  "
    "D30" = "0.131"
    "E30" = "  +11.38%  "
    "E31" = "  -1.91%  "
    "E32" = "  -1.09%  "
    "D33" = "4.92"
    "E33" = "  -7.13%  "
    "D34" = "0.0630"
    "E34" = "  +0.30%  "
    "E35" = "  -6.26%  "
    "D36" = "6.33"
    "E36" = "  +4.56%  "
    "E37" = "  -0.08%  "
    "E38" = "  -7.59%  "
    "E39" = "  -3.14%  "
    "D40" = "3.04"
    "E40" = "  +2.08%  "
    "E41" = "  -0.49%  "
    "D42" = "0.0961"
    "E42" = "  -8.06%  "
    "D44" = "0.0214"
    "E44" = "  -2.32%  "
    "E45" = "  -5.07%  "
    "D46" = "92.81"
    "E46" = "  -4.29%  "
    "D47" = "16.23"
    "E47" = "  -5.85%  "
    "D48" = "1.382.71"
    "E48" = "  -2.74%  "
    "D49" = "7.51"
    "E49" = "  -6.07%  "
    "E50" = "  -3.27%  "
    "D51" = "45.51"
    "E51" = "  -2.88%  "
}

foreach ($cell in $updates.Keys) {
    if ($textCells -contains $cell) {
        $ws.Range($cell).NumberFormat = "@"
        $ws.Range($cell).Value = $updates[$cell]
        $ws.Range($cell).Style = "Normal"
    } else {
        $ws.Range($cell).Value = $updates[$cell]
    }
}
